$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search")

# Row 17 (S1_TC_T16): change sort order from score:asc -> score:desc, and fix status to PASS
$ws.Range("B17").Value = "Search for documents and sort on score - desc"
$ws.Range("G17").Value = "?query=biotechnology&sort=_score:desc"
$ws.Range("L17").Value = "PASS"

# Row 18 (S1_TC_T17): shift content that used to be "score - desc" to "pub date - des(c)"
$ws.Range("B18").Value = "Search for documents and sort on pub date - des"
$ws.Range("G18").Value = "?query=biotechnology&sort=sortdate:desc"

# Row 19 (S1_TC_T18): shift content that used to be row 20's "pub date - asc"
$ws.Range("B19").Value = "Search for documents and sort on pub date - asc"
$ws.Range("G19").Value = "?query=biotechnology&sort=sortdate:asc"

# Row 20 (old S1_TC_T19) is now redundant and gets removed entirely
$ws.Rows.Item(20).Delete()

# Update sheet view / selection to match target state
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("A2:A19").Select()
